# Gráficos para combinaciones de ángulos.
# Update the min/max run id values on Sheet1 and move the active
# selection from D3 to D2 to match the authored change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Id_max_run (C2) and the corresponding D2 value.
$ws.Range("C2").Value = 547
$ws.Range("D2").Value = 550

# Make sure Sheet1 is the active sheet and move the selection to D2.
$ws.Activate()
$ws.Range("D2").Select() | Out-Null
